$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 491, shifting existing rows 491:524 down to 494:527
$ws.Range("A491:R493").EntireRow.Insert()

# New row 491 - Melón / Calameño / Primera - Provincia de Quillota
$ws.Cells.Item(491,1).Value2  = 3
$ws.Cells.Item(491,2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(491,3).Value2  = "Coquimbo"
$ws.Cells.Item(491,4).Value2  = 44931
$ws.Cells.Item(491,5).Value2  = 5
$ws.Cells.Item(491,6).Value2  = 100112027
$ws.Cells.Item(491,7).Value2  = "Melón"
$ws.Cells.Item(491,8).Value2  = "Calameño"
$ws.Cells.Item(491,9).Value2  = "Primera"
$ws.Cells.Item(491,10).Value2 = 400
$ws.Cells.Item(491,11).Value2 = 1800
$ws.Cells.Item(491,12).Value2 = 1800
$ws.Cells.Item(491,13).Value2 = 1800
$ws.Cells.Item(491,14).Value2 = "$/unidad"
$ws.Cells.Item(491,15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(491,16).Value2 = 1800
$ws.Cells.Item(491,17).Value2 = 1
$ws.Cells.Item(491,18).Value2 = "Hortaliza"

# New row 492 - Melón / Calameño / Primera - Provincia de Quillota
$ws.Cells.Item(492,1).Value2  = 3
$ws.Cells.Item(492,2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(492,3).Value2  = "Coquimbo"
$ws.Cells.Item(492,4).Value2  = 44931
$ws.Cells.Item(492,5).Value2  = 5
$ws.Cells.Item(492,6).Value2  = 100112027
$ws.Cells.Item(492,7).Value2  = "Melón"
$ws.Cells.Item(492,8).Value2  = "Calameño"
$ws.Cells.Item(492,9).Value2  = "Primera"
$ws.Cells.Item(492,10).Value2 = 1030
$ws.Cells.Item(492,11).Value2 = 1500
$ws.Cells.Item(492,12).Value2 = 1700
$ws.Cells.Item(492,13).Value2 = 1593
$ws.Cells.Item(492,14).Value2 = "$/unidad"
$ws.Cells.Item(492,15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(492,16).Value2 = 1593
$ws.Cells.Item(492,17).Value2 = 1
$ws.Cells.Item(492,18).Value2 = "Hortaliza"

# New row 493 - Melón / Calameño / Segunda - Provincia de Quillota
$ws.Cells.Item(493,1).Value2  = 3
$ws.Cells.Item(493,2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(493,3).Value2  = "Coquimbo"
$ws.Cells.Item(493,4).Value2  = 44931
$ws.Cells.Item(493,5).Value2  = 5
$ws.Cells.Item(493,6).Value2  = 100112027
$ws.Cells.Item(493,7).Value2  = "Melón"
$ws.Cells.Item(493,8).Value2  = "Calameño"
$ws.Cells.Item(493,9).Value2  = "Segunda"
$ws.Cells.Item(493,10).Value2 = 530
$ws.Cells.Item(493,11).Value2 = 1000
$ws.Cells.Item(493,12).Value2 = 1000
$ws.Cells.Item(493,13).Value2 = 1000
$ws.Cells.Item(493,14).Value2 = "$/unidad"
$ws.Cells.Item(493,15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(493,16).Value2 = 1000
$ws.Cells.Item(493,17).Value2 = 1
$ws.Cells.Item(493,18).Value2 = "Hortaliza"

Write-Output "inserted 3 rows and populated rows 491-493"
